$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column A width to match the target layout (closest achievable value
# given the engine's internal width quantization)
$ws.Columns.Item(1).ColumnWidth = 67.3

# Add the "Checked" header in B1
$ws.Range("B1").Value = "Checked"

# Fill B2:B13 with 0
$ws.Range("B2:B13").Value = 0

# Update selection to match the target view state
$ws.Range("D15").Select()
